# Program Edit feature: update existing Program records on the "Program" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Program")

# Row 6 / Row 7 together represent a second program record whose name and
# description are being corrected.
$ws.Range("A6").Value = "ABC345"
$ws.Range("B7").Value = "ABC Description"

# Row 5 holds the first program record; its ProgramName value is updated.
$ws.Range("A5").Value = "SM 000001235645"
